$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BO2").Value = 0.7837604284286499
$ws.Range("BO3").Value = 0.7279143333435059
$ws.Range("BO4").Value = 0.813188374042511
$ws.Range("BO5").Value = 0.876110315322876
$ws.Range("BO6").Value = 0.741787314414978
$ws.Range("BO7").Value = 0.7662057876586914
$ws.Range("BO8").Value = 0.7803550958633423
$ws.Range("BO9").Value = 0.7516946196556091
$ws.Range("BO10").Value = 0.7472559809684753
$ws.Range("BO11").Value = 0.7001696825027466
$ws.Range("BO12").Value = 0.7348120212554932
$ws.Range("BO13").Value = 0.7973758578300476
$ws.Range("BO14").Value = 0.7539451122283936
$ws.Range("BO15").Value = 0.7015854716300964
$ws.Range("BO16").Value = 0.7444555759429932
$ws.Range("BO17").Value = 0.7379483580589294
$ws.Range("BO18").Value = 0.7753097414970398
$ws.Range("BO19").Value = 0.4195190370082855
$ws.Range("BO20").Value = 0.8011788129806519
$ws.Range("BO21").Value = 0.456044465303421
$ws.Range("BO22").Value = 0.5989366769790649
$ws.Range("BO23").Value = 0.7619568705558777
$ws.Range("BO24").Value = 0.7488059401512146
$ws.Range("BO25").Value = 0.6251933574676514
$ws.Range("BO26").Value = 0.7168236970901489
$ws.Range("BO27").Value = 0.64058518409729
$ws.Range("BO28").Value = 0.8037018179893494
$ws.Range("BO29").Value = 0.625207245349884
$ws.Range("BO30").Value = 0.8132553100585938
$ws.Range("BO31").Value = 0.7647276520729065
$ws.Range("BO32").Value = 0.7053089737892151
$ws.Range("BO33").Value = 0.7426403760910034
$ws.Range("BO34").Value = 0.592975378036499
$ws.Range("BO35").Value = 0.7315511107444763
$ws.Range("BO36").Value = 0.5152532458305359
$ws.Range("BO37").Value = 0.7740731835365295
$ws.Range("BO38").Value = 0.7825538516044617
$ws.Range("BO39").Value = 0.462871640920639
$ws.Range("BO40").Value = 0.6657180786132812
$ws.Range("BO41").Value = 0.5833063125610352
$ws.Range("BO42").Value = 0.4351317882537842
$ws.Range("BO43").Value = 0.5796887874603271
$ws.Range("BO44").Value = 0.6429891586303711
$ws.Range("BO45").Value = 0.5113218426704407
$ws.Range("BO46").Value = 0.5236095190048218
$ws.Range("BO47").Value = 0.8046386241912842
$ws.Range("BO48").Value = 0.8060837388038635
$ws.Range("BO49").Value = 0.7836665511131287
$ws.Range("BO50").Value = 0.7159928679466248
$ws.Range("BO51").Value = 0.7895399928092957
$ws.Range("BO52").Value = 0.4180402159690857
$ws.Range("BO53").Value = 0.2976245582103729
$ws.Range("BO54").Value = 0.01408990658819675
$ws.Range("BO55").Value = 0.7957469820976257
$ws.Range("BO56").Value = 0.8261958956718445
$ws.Range("BO57").Value = 0.8513769507408142
$ws.Range("BO58").Value = 0.8695443868637085
$ws.Range("BO59").Value = 0.8374533653259277
$ws.Range("BO60").Value = 0.7939706444740295
$ws.Range("BO61").Value = 0.8358448147773743
$ws.Range("BO62").Value = 0.7842193245887756
$ws.Range("BO63").Value = 0.7223101258277893
$ws.Range("BO64").Value = 0.5808306932449341
$ws.Range("BO65").Value = 0.6600769758224487
$ws.Range("BO66").Value = 0.145443394780159
$ws.Range("BO67").Value = 0.5602331757545471
$ws.Range("BO68").Value = 0.1382258832454681
$ws.Range("BO69").Value = 0.7583444714546204
$ws.Range("BO70").Value = 0.6827608942985535
$ws.Range("BO71").Value = 0.3233097195625305
$ws.Range("BO72").Value = 0.7946881651878357
$ws.Range("BO73").Value = 0.7534769177436829
$ws.Range("BO74").Value = 0.7075704336166382
$ws.Range("BO75").Value = 0.8701086640357971
$ws.Range("BO76").Value = 0.821523129940033
